$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price column (D) keeps its original text formatting so that
# numeric-looking price strings (e.g. "576.99") are not auto-converted to
# real numbers when assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.754.12"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "3.151.45"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "576.99"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "148.84"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.151.15"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("D11").Value = "6.09"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").Value = "0.498"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").Value = "36.96"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "3.663.22"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "64.898.27"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "3.148.58"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "7.09"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "502.45"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("D21").Value = "14.79"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "0.711"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").Value = "7.69"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").Value = "83.85"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").Value = "8.86"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("E30").Value = "  +5.83%  "
$ws.Range("D31").Value = "27.44"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "1.19"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("D34").Value = "6.13"
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("D35").Value = "6.44"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("D36").Value = "54.55"
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("D37").Value = "0.0890"
$ws.Range("E37").Value = "  +3.59%  "
$ws.Range("D38").Value = "473.85"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").Value = "0.0414"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("D40").Value = "2.92"
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("D41").Value = "8.63"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "3.003.42"
$ws.Range("E42").Value = "  -3.63%  "
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("D44").Value = "2.41"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("D45").Value = "0.280"
$ws.Range("E45").Value = "  -3.97%  "
$ws.Range("D46").Value = "28.09"
$ws.Range("E46").Value = "  -3.39%  "
$ws.Range("D47").Value = "0.0₃0581"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "0.114"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("E50").Value = "  -3.11%  "
$ws.Range("D51").Value = "33.77"
$ws.Range("E51").Value = "  +7.90%  "
